# Update "eCoachingLog Lessons Learned" workbook with the new lesson-learned
# entry (added 5/14/2015) and retitle the sheet heading.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new "Maintenance" phase lesson-learned row (row 4) ---------
# Fill these in first so that the new shared strings get appended to the
# shared-string table in the same order they end up in the target file
# (row 4 contents, then the retitled heading last).
$ws.Range("A4").Value = "Regression testing should be performed when major changes are made"
$ws.Range("B4").Value = "2015-05-14"
$ws.Range("C4").Value = "Tim Queen"
$ws.Range("D4").Value = "SPM"
$ws.Range("E4").Value = "Changes made to one part of the system may have a negative impact on other parts of the system"
$ws.Range("F4").Value = "Create regression test cases that can be used when significant changes are made to eCL."
$ws.Range("G4").Value = "Maintenance"

# --- Rename the report heading from "Performance Management Lessons
#     Learned" to "eCoaching Log Lessons Learned" --------------------------
$ws.Range("A1").Value = "eCoaching Log Lessons Learned"

# --- Leave the selection on G2, matching the author's last click ----------
$ws.Range("G2").Select() | Out-Null
